# Update "想去人数" (F column) values across sheets to reflect refreshed
# generated output, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 24
$ws1.Range("F5").Value = 15317
$ws1.Range("F6").Value = 413
$ws1.Range("F9").Value = 15323
$ws1.Range("F10").Value = 46
$ws1.Range("F11").Value = 8898
$ws1.Range("F15").Value = 191
$ws1.Range("F17").Value = 190
$ws1.Range("F19").Value = 40
$ws1.Range("F20").Value = 533
$ws1.Range("F22").Value = 9
$ws1.Range("F23").Value = 57
$ws1.Range("F24").Value = 1101
$ws1.Range("F26").Value = 18
$ws1.Range("F27").Value = 67
$ws1.Range("F29").Value = 36
$ws1.Range("F31").Value = 41
$ws1.Range("F34").Value = 292
$ws1.Range("F36").Value = 112
$ws1.Range("F37").Value = 5439

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 65

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 24
$ws4.Range("F5").Value = 15317
$ws4.Range("F6").Value = 413
$ws4.Range("F9").Value = 15323
$ws4.Range("F10").Value = 46
$ws4.Range("F11").Value = 8898
$ws4.Range("F16").Value = 191
$ws4.Range("F18").Value = 190
$ws4.Range("F20").Value = 40
$ws4.Range("F21").Value = 533
$ws4.Range("F23").Value = 9
$ws4.Range("F24").Value = 57
$ws4.Range("F25").Value = 1101
$ws4.Range("F27").Value = 18
$ws4.Range("F28").Value = 67
$ws4.Range("F30").Value = 36
$ws4.Range("F31").Value = 65
$ws4.Range("F34").Value = 41
$ws4.Range("F37").Value = 292
$ws4.Range("F39").Value = 112
$ws4.Range("F40").Value = 5439
